$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21: Book and a Hard Place / Engraved Hard Leather Grimoire
$ws.Range("H21").Value = 9900
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 9900
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 9900
$ws.Range("N21").Value = -10836
$ws.Range("M21").ClearContents()

# Row 23: There's Something about Bury / Hard Leather Grimoire
$ws.Range("H23").Value = 9900
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 9900
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 9900
$ws.Range("N23").Value = -10368
$ws.Range("M23").ClearContents()

# Row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Range("H38").Value = 835.0909
$ws.Range("I38").Value = 76.22221999999999
$ws.Range("J38").Value = 4250
$ws.Range("K38").Value = 228.66666
$ws.Range("L38").Value = 12750
$ws.Range("M38").Value = 143.33334
$ws.Range("N38").Value = -13494

# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 1423.25
$ws.Range("I43").Value = 3333
$ws.Range("J43").Value = 1150.4286
$ws.Range("K43").Value = 3333
$ws.Range("L43").Value = 1150.4286
$ws.Range("M43").Value = -3264
$ws.Range("N43").Value = -1288.4286

# Row 44: Alive and Unwell / Budding Oak Wand
$ws.Range("H44").Value = 19999.25
$ws.Range("J44").Value = 19999.25
$ws.Range("L44").Value = 19999.25
$ws.Range("N44").Value = -20923.25

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 12347179

# Row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 1420.7
$ws.Range("I58").Value = 452.16666
$ws.Range("J58").Value = 2873.5
$ws.Range("K58").Value = 1356.49998
$ws.Range("L58").Value = 8620.5
$ws.Range("M58").Value = -1206.49998
$ws.Range("N58").Value = -8920.5

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 2623.7727
$ws.Range("I69").Value = 2082.6
$ws.Range("J69").Value = 2782.9412
$ws.Range("K69").Value = 6247.799999999999
$ws.Range("L69").Value = 8348.8236
$ws.Range("M69").Value = -5373.799999999999
$ws.Range("N69").Value = -10096.8236

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 2623.7727
$ws.Range("I72").Value = 2082.6
$ws.Range("J72").Value = 2782.9412
$ws.Range("K72").Value = 18743.4
$ws.Range("L72").Value = 25046.4708
$ws.Range("M72").Value = -14375.4
$ws.Range("N72").Value = -33782.4708

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 1304.098
$ws.Range("I80").Value = 1027.2572
$ws.Range("J80").Value = 1909.6875
$ws.Range("K80").Value = 3081.7716
$ws.Range("L80").Value = 5729.0625
$ws.Range("M80").Value = -2083.7716
$ws.Range("N80").Value = -7725.0625

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 1304.098
$ws.Range("I83").Value = 1027.2572
$ws.Range("J83").Value = 1909.6875
$ws.Range("K83").Value = 9245.3148
$ws.Range("L83").Value = 17187.1875
$ws.Range("M83").Value = -4253.3148
$ws.Range("N83").Value = -27171.1875

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 10563334
$ws.Range("I86").Value = 12072239
$ws.Range("J86").Value = 1002
$ws.Range("K86").Value = 12072239
$ws.Range("L86").Value = 1002
$ws.Range("M86").Value = -12071116
$ws.Range("N86").Value = -3248

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 10563334
$ws.Range("I89").Value = 12072239
$ws.Range("J89").Value = 1002
$ws.Range("K89").Value = 60361195
$ws.Range("L89").Value = 5010
$ws.Range("M89").Value = -60355579
$ws.Range("N89").Value = -16242

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 6441344.5
$ws.Range("I116").Value = 7085244
$ws.Range("J116").Value = 2350
$ws.Range("K116").Value = 7085244
$ws.Range("L116").Value = 2350
$ws.Range("M116").Value = -7081802
$ws.Range("N116").Value = -9234

$ws = $wb.Worksheets.Item("ARM")
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 7938857
$ws.Range("I97").Value = 8335749.5
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 8335749.5
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -8335253.5
$ws.Range("N97").Value = -1992

$ws = $wb.Worksheets.Item("BSM")
# Row 40: Can You Spare a Dolabra / Steel Dolabra
$ws.Range("H40").Value = 28870.084
$ws.Range("J40").Value = 28870.084
$ws.Range("L40").Value = 28870.084
$ws.Range("N40").Value = -29400.084

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3527.625
$ws.Range("I86").Value = 3700
$ws.Range("J86").Value = 3355.25
$ws.Range("K86").Value = 3700
$ws.Range("L86").Value = 3355.25
$ws.Range("M86").Value = -2577
$ws.Range("N86").Value = -5601.25

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3527.625
$ws.Range("I89").Value = 3700
$ws.Range("J89").Value = 3355.25
$ws.Range("K89").Value = 18500
$ws.Range("L89").Value = 16776.25
$ws.Range("M89").Value = -12884
$ws.Range("N89").Value = -28008.25

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 8980.416999999999
$ws.Range("I94").Value = 683.25
$ws.Range("J94").Value = 50466.25
$ws.Range("K94").Value = 683.25
$ws.Range("L94").Value = 50466.25
$ws.Range("M94").Value = -232.25
$ws.Range("N94").Value = -51368.25

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1618.3334
$ws.Range("I99").Value = 1385.5555
$ws.Range("J99").Value = 1851.1111
$ws.Range("K99").Value = 1385.5555
$ws.Range("L99").Value = 1851.1111
$ws.Range("M99").Value = 112.4445000000001
$ws.Range("N99").Value = -4847.1111

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 2394.6155
$ws.Range("I107").Value = 1907.1428
$ws.Range("J107").Value = 2963.3333
$ws.Range("K107").Value = 1907.1428
$ws.Range("L107").Value = 2963.3333
$ws.Range("M107").Value = 12.85719999999992
$ws.Range("N107").Value = -6803.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 23811530
$ws.Range("I16").Value = 33335140
$ws.Range("J16").Value = 2504.3333
$ws.Range("K16").Value = 33335140
$ws.Range("L16").Value = 2504.3333
$ws.Range("M16").Value = -33334853
$ws.Range("N16").Value = -3078.3333

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 23811530
$ws.Range("I113").Value = 33335140
$ws.Range("J113").Value = 2504.3333
$ws.Range("K113").Value = 33335140
$ws.Range("L113").Value = 2504.3333
$ws.Range("M113").Value = -33332970
$ws.Range("N113").Value = -6844.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 872.5700000000001
$ws.Range("J131").Value = 948.3837
$ws.Range("L131").Value = 2845.1511
$ws.Range("N131").Value = -12925.1511

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3535.88
$ws.Range("I80").Value = 4335.909
$ws.Range("J80").Value = 2907.2856
$ws.Range("K80").Value = 4335.909
$ws.Range("L80").Value = 2907.2856
$ws.Range("M80").Value = -3337.909
$ws.Range("N80").Value = -4903.2856

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3535.88
$ws.Range("I83").Value = 4335.909
$ws.Range("J83").Value = 2907.2856
$ws.Range("K83").Value = 21679.545
$ws.Range("L83").Value = 14536.428
$ws.Range("M83").Value = -16687.545
$ws.Range("N83").Value = -24520.428

$ws = $wb.Worksheets.Item("LTW")
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 1003400.4
$ws.Range("I68").Value = 1669000.6
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 1669000.6
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -1668251.6
$ws.Range("N68").Value = -6498

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 1003400.4
$ws.Range("I71").Value = 1669000.6
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 8345003
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -8341259
$ws.Range("N71").Value = -32488

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1500.8966
$ws.Range("I93").Value = 1405.25
$ws.Range("J93").Value = 1960
$ws.Range("K93").Value = 1405.25
$ws.Range("L93").Value = 1960
$ws.Range("M93").Value = -157.25
$ws.Range("N93").Value = -4456

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 6464.4614
$ws.Range("I136").Value = 9556.799999999999
$ws.Range("J136").Value = 2247.6365
$ws.Range("K136").Value = 28670.4
$ws.Range("L136").Value = 6742.9095
$ws.Range("M136").Value = -26120.4
$ws.Range("N136").Value = -11842.9095

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 2199
$ws.Range("I81").Value = 1982.25
$ws.Range("K81").Value = 3964.5
$ws.Range("M81").Value = -2903.5

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 2199
$ws.Range("I84").Value = 1982.25
$ws.Range("K84").Value = 19822.5
$ws.Range("M84").Value = -14518.5

# Row 86: Felt for the Fallen / Chimerical Felt
$ws.Range("H86").Value = 4162.5
$ws.Range("J86").Value = 4162.5
$ws.Range("L86").Value = 4162.5
$ws.Range("N86").Value = -6408.5

# Row 89: Blinded Veil of Vigilance (L) / Chimerical Felt
$ws.Range("H89").Value = 4162.5
$ws.Range("J89").Value = 4162.5
$ws.Range("L89").Value = 20812.5
$ws.Range("N89").Value = -32044.5

# Row 94: Proper Props / Bloodhempen Armguards of Scouting
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 715.88464
$ws.Range("I107").Value = 690.36365
$ws.Range("J107").Value = 734.6
$ws.Range("K107").Value = 2071.09095
$ws.Range("L107").Value = 2203.8
$ws.Range("M107").Value = -151.0909499999998
$ws.Range("N107").Value = -6043.8
